# ErrorMachineTemplate.xlsx - adjust import/export excel attribute names.
#
# The hidden header row (row 2) carries the machine-readable field names used
# by the import/export mapping. Column B's field name is being renamed from
# the generic "name" to the more specific "ten_su_co" (Vietnamese for
# "error name") to line up with the other already-specific field names in
# that row (id, line_name, nguyen_nhan, cach_xu_ly).
#
# The visible header row (row 4) keeps its same Vietnamese captions
# (Mã lỗi, Tên lỗi, Công đoạn, Nguyên nhân, Cách xử lý) - only the shared
# string table is reshuffled as a side effect of removing the no-longer-used
# "name" string and appending "ten_su_co".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the hidden machine field name in column B from "name" to "ten_su_co".
$ws.Range("B2").Value = "ten_su_co"

# Update the active selection/cursor position left in the sheet.
[void]$ws.Range("E12").Select()
